# Insert a new row above row 274, duplicating the row that is currently
# there (it will be pushed down to row 275 by the insert) and giving the
# new row a fresh "Fecha" (column D) value. Every row from the old 274
# through the old 415 shifts down by one, ending at new row 416.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("274:274").Insert()

$lastCol = 18
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item(274, $c).Value2 = $ws.Cells.Item(275, $c).Value2
}

# New row's date (column D / 4): serial 45029 = 2023-04-13
$ws.Cells.Item(274, 4).Value2 = 45029
